$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.823.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.811.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4304'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3685'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07201'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8660'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.86'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.948.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.621'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.341'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06890'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008888'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.838.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.207'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.159.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.867'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.208'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.900'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08912'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7575'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.73%  '

$ws.Range("E33").Value = '  +6.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.427'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.768'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.005'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.122'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05211'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5075'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1644'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.655'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.523'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.280'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.36%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.17%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.653'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.87%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06269'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4542'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.806'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.32%  '
